$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (row 2)
$ws.Range("G2").Value = 5.106677666666667
$ws.Range("H2").Value = 15.320033
$ws.Range("I2").Value = 0.2427554644382946
$ws.Range("J2").Value = 0.2427554644382946
$ws.Range("M2").Value = 0.3947560000000001
$ws.Range("N2").Value = 1.184268
$ws.Range("O2").Value = 0.2067464753251374
$ws.Range("P2").Value = 0.2067464753251374
$ws.Range("Q2").Value = 2.015891648982667
$ws.Range("R2").Value = 18.143024840844
$ws.Range("S2").Value = 0.05018883663853415
$ws.Range("T2").Value = 0.05018883663853414

# row 3
$ws.Range("G3").Value = 5.106677666666667
$ws.Range("H3").Value = 15.320033
$ws.Range("I3").Value = 0.2427554644382946
$ws.Range("J3").Value = 0.2427554644382946
$ws.Range("O3").Value = 0.676395925572051
$ws.Range("P3").Value = 0.676395925572051
$ws.Range("Q3").Value = 6.595231650852779
$ws.Range("R3").Value = 59.35708485767501
$ws.Range("S3").Value = 0.1641988070564134
$ws.Range("T3").Value = 0.1641988070564134

# row 4
$ws.Range("G4").Value = 5.106677666666667
$ws.Range("H4").Value = 15.320033
$ws.Range("I4").Value = 0.2427554644382946
$ws.Range("J4").Value = 0.2427554644382946
$ws.Range("M4").Value = 0.2231246666666666
$ws.Range("N4").Value = 0.6693739999999999
$ws.Range("O4").Value = 0.1168575991028116
$ws.Range("P4").Value = 0.1168575991028116
$ws.Range("Q4").Value = 1.139425752149111
$ws.Range("R4").Value = 10.254831769342
$ws.Range("S4").Value = 0.02836782074334707
$ws.Range("T4").Value = 0.02836782074334707

# row 5
$ws.Range("I5").Value = 0.6077416388947038
$ws.Range("J5").Value = 0.6077416388947038
$ws.Range("M5").Value = 0.3947560000000001
$ws.Range("N5").Value = 1.184268
$ws.Range("O5").Value = 0.2067464753251374
$ws.Range("P5").Value = 0.2067464753251374
$ws.Range("Q5").Value = 5.046812426742667
$ws.Range("R5").Value = 45.42131184068401
$ws.Range("S5").Value = 0.1256484417498024
$ws.Range("T5").Value = 0.1256484417498024

# row 6
$ws.Range("I6").Value = 0.6077416388947038
$ws.Range("J6").Value = 0.6077416388947038
$ws.Range("O6").Value = 0.676395925572051
$ws.Range("P6").Value = 0.676395925572051
$ws.Range("S6").Value = 0.4110739683488583
$ws.Range("T6").Value = 0.4110739683488583

# row 7
$ws.Range("I7").Value = 0.6077416388947038
$ws.Range("J7").Value = 0.6077416388947038
$ws.Range("M7").Value = 0.2231246666666666
$ws.Range("N7").Value = 0.6693739999999999
$ws.Range("O7").Value = 0.1168575991028116
$ws.Range("P7").Value = 0.1168575991028116
$ws.Range("S7").Value = 0.07101922879604299
$ws.Range("T7").Value = 0.07101922879604299

# row 8
$ws.Range("I8").Value = 0.1495028966670016
$ws.Range("J8").Value = 0.1495028966670016
$ws.Range("M8").Value = 0.3947560000000001
$ws.Range("N8").Value = 1.184268
$ws.Range("O8").Value = 0.2067464753251374
$ws.Range("P8").Value = 0.2067464753251374
$ws.Range("Q8").Value = 1.241503014513334
$ws.Range("R8").Value = 11.17352713062
$ws.Range("S8").Value = 0.03090919693680081
$ws.Range("T8").Value = 0.03090919693680081

# row 9
$ws.Range("I9").Value = 0.1495028966670016
$ws.Range("J9").Value = 0.1495028966670016
$ws.Range("O9").Value = 0.676395925572051
$ws.Range("P9").Value = 0.676395925572051
$ws.Range("S9").Value = 0.1011231501667792
$ws.Range("T9").Value = 0.1011231501667792

# row 10
$ws.Range("I10").Value = 0.1495028966670016
$ws.Range("J10").Value = 0.1495028966670016
$ws.Range("M10").Value = 0.2231246666666666
$ws.Range("N10").Value = 0.6693739999999999
$ws.Range("O10").Value = 0.1168575991028116
$ws.Range("P10").Value = 0.1168575991028116
$ws.Range("Q10").Value = 0.7017244735455554
$ws.Range("R10").Value = 6.315520261909999
$ws.Range("S10").Value = 0.01747054956342154
$ws.Range("T10").Value = 0.01747054956342154
